$d = $word.ActiveDocument
$found = $d.Content.Find.Execute("tutorials by means of sponsorships", $true, $false, $false, $false, $false, $true, 1, $false, "events by means of sponsorships", 2)
Write-Host "Found:" $found
